# Add a new "affix_type" column (column AN / 40) to the Affixes sheet and
# populate it with the affix type filter value (7) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Affixes")

# Insert a brand new column at AN (40) right after the existing data (AM/39).
$ws.Columns.Item(40).Insert()

$headerRange = $ws.Range("AN1")
$dataRange = $ws.Range("AN2:AN10")

# Match the workbook's default "Normal" cell style (Calibri 11 / General)
# instead of inheriting the Arial-based style of the neighbouring column.
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 11
$headerRange.NumberFormat = "General"

$dataRange.Font.Name = "Calibri"
$dataRange.Font.Size = 11
$dataRange.NumberFormat = "General"

# New header label.
$headerRange.Value = "affix_type"

# Every existing affix row gets the same affix_type value (7).
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 40).Value = 7
}

# Move the active selection to the newly added column/row, matching the
# author's final cursor position.
$ws.Range("AN10").Select()
